$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "policy search" column for the SIS (lattice, L=100) table ---
$ws.Range("G3").Value = "0.03(0.01)"
$ws.Range("G4").Value = "0.40(0.02)"
$ws.Range("G5").Value = "0.43(0.01)"

# --- New "Ebola" table at rows 15-17 ---
$ws.Range("A15").Value = "Ebola"
$ws.Range("A15").Font.Bold = $true

$ws.Range("C17").Value = "0.17(0.01)"
$ws.Range("D17").Value = "0.21(0.005)"
$ws.Range("G17").Value = "0.16(0.004)"
$ws.Range("B17").Value = "0.36(0.003)"

$ws.Range("B16").Value = "random"
$ws.Range("C16").Value = "mf one step"
$ws.Range("D16").Value = "mb one step"
$ws.Range("E16").Value = "mse-averaged one step"
$ws.Range("F16").Value = "equal-averaged one step"
$ws.Range("G16").Value = "policy search"

# --- Update "policy search" column value in last row, and selection ---
$ws.Range("G6").Value = "0.45(0.02)"

# --- Update active cell selection to G6 ---
$ws.Range("G6").Select() | Out-Null
